# Generate Report for Handback
# Update the timestamp strings recorded for the 74e4990e-... handback row
# across the Overview, zh-cn, and de-de sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for 74e4990e... row (row 3)
$overview.Range("G3").Value = "2016-09-05 14:57:22"

# zh-cn sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime" for 74e4990e... row (row 3)
$zhcn.Range("H3").Value = "2016-09-05 14:57:17"
$zhcn.Range("K3").Value = "2016-09-05 14:57:35"

# de-de sheet: "Correspond Handoff Datetime" / "Correspond Handback DateTime" for 74e4990e... row (row 3)
$dede.Range("H3").Value = "2016-09-05 14:57:22"
$dede.Range("K3").Value = "2016-09-05 14:57:43"
